# Applies the betexplorer "slovenia 2-snl 2023-2024" update:
#  - Several existing match rows had their F:V payload (home..url_partida)
#    rotated among the rows of the same fixture-date block - the rows were
#    mis-assigned to the wrong fixture and this corrects it. Columns A:E
#    (Indice/pais/torneio/temporada/data_partida) stay on their own row.
#  - One brand-new match row (138) is appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-Rows($RowCycle) {
    # For cycle (a, b, c, ...): new_a = old_b, new_b = old_c, ..., new_last = old_a
    $n = $RowCycle.Length
    $first = $RowCycle[0]
    $saved = $ws.Range("F" + $first + ":V" + $first).Value2

    for ($i = 0; $i -lt ($n - 1); $i++) {
        $dstRow = $RowCycle[$i]
        $srcRow = $RowCycle[$i + 1]
        $srcVals = $ws.Range("F" + $srcRow + ":V" + $srcRow).Value2
        $ws.Range("F" + $dstRow + ":V" + $dstRow).Value2 = $srcVals
    }

    $lastRow = $RowCycle[$n - 1]
    $ws.Range("F" + $lastRow + ":V" + $lastRow).Value2 = $saved
}

# Row-content rotations (derived from the canonical-XML diff)
Rotate-Rows @(67, 68)
Rotate-Rows @(69, 70, 71)
Rotate-Rows @(73, 74)
Rotate-Rows @(76, 78, 79)
Rotate-Rows @(91, 93)
Rotate-Rows @(92, 94)
Rotate-Rows @(95, 96)
Rotate-Rows @(113, 114, 115)
Rotate-Rows @(117, 120)
Rotate-Rows @(130, 133, 131, 132)
Rotate-Rows @(134, 136, 135)

# Append the new match row 138, copying number formats/styles from row 137
# so the new "Indice" (A) and "data_partida" (E) cells match the rest of
# the table's formatting.
$ws.Range("A137").Copy() | Out-Null
$ws.Range("A138").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E137").Copy() | Out-Null
$ws.Range("E138").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A138").Value = 137
$ws.Range("B138").Value = "slovenia"
$ws.Range("C138").Value = "2-snl"
$ws.Range("D138").Value = "2023-2024"
$ws.Range("E138").Value = 45247.625
$ws.Range("F138").Value = "NK Krka"
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = "Dravinja"
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = 1.72
$ws.Range("K138").Value = "16/11/2023 03:12"
$ws.Range("L138").Value = 2.34
$ws.Range("M138").Value = "17/11/2023 14:44"
$ws.Range("N138").Value = 3.59
$ws.Range("O138").Value = "16/11/2023 03:12"
$ws.Range("P138").Value = 3.54
$ws.Range("Q138").Value = "17/11/2023 14:46"
$ws.Range("R138").Value = 3.71
$ws.Range("S138").Value = "16/11/2023 03:12"
$ws.Range("T138").Value = 2.63
$ws.Range("U138").Value = "17/11/2023 14:44"
$ws.Range("V138").Value = "https://www.betexplorer.com/football/slovenia/2-snl/nk-krka-dravinja/Ctod0zgQ/"
